# IPA.xlsx - "Checking in the my scripts"
# Adds a new IPA012 test-case row (NEON-291/400/438/574) and flips the
# Runmode flag on the two existing rows from Y -> N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- New row 4: IPA012 test case ---------------------------------------
$ws.Range("B4").Value = "NEON-291||NEON-400||NEON-438||NEON-574"
$ws.Range("C4").Value = "User must be able to form a Technology Search||User must be able to return to the app landing page via the App header to start a new search||Option on the IPA App header to allow the user to return to the app landing page||User must be able to form a Company Search"
$ws.Range("A4").Value = "IPA012"

# --- Existing rows: Runmode (col D) Y -> N -----------------------------
$ws.Range("D2").Value = "N"
$ws.Range("D3").Value = "N"

# --- Finish filling in the new row --------------------------------------
$ws.Range("E4").Value = "FAIL"
$ws.Range("D4").Value = "Y"

# Match row 4 formatting (thin border, no fill) to the row above it
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths (best-fit style, now that col B/E hold new text) ----
$ws.Columns.Item(2).ColumnWidth = 44.140625
$ws.Columns.Item(5).ColumnWidth = 7.2578125

# --- Selection as left by the editing session ---------------------------
$ws.Range("D2").Select()
